# Rellenar la columna "uds. Objetivo semana pasada" (columna Q) con los
# valores reales de la semana pasada para las filas que tenian 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_8")

$valores = @{
    2  = 2
    3  = 103
    4  = 12
    5  = 4
    8  = 5
    9  = 2
    11 = 3
    12 = 4
    13 = 2
    14 = 4
    19 = 3
    21 = 3
    22 = 10
    23 = 8
    24 = 2
    25 = 4
    28 = 17
    29 = 8
    32 = 24
    34 = 2
    35 = 2
}

foreach ($fila in $valores.Keys) {
    $ws.Range("Q$fila").Value = $valores[$fila]
}
